$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data, and the Polkadot/BabyDogeCoin row swap.

$ws.Range('D2').Value = '62.745.77'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').Value = '2.444.27'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'575.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').Value = "'145.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '2.442.67'
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('D11').Value = "'0.163"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.24%  '
$ws.Range('D12').Value = "'5.24"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').Value = "'28.15"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.98%  '
$ws.Range('D15').Value = "'0.0000179"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.22%  '
$ws.Range('D16').Value = '2.889.06'
$ws.Range('E16').Value = '  +4.01%  '
$ws.Range('D17').Value = '62.710.19'
$ws.Range('E17').Value = '  +3.92%  '
$ws.Range('D18').Value = '2.446.62'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').Value = "'7.85"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('D20').Value = "'10.98"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.89%  '
$ws.Range('D21').Value = "'329.65"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = "'4.14"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('B23').Value = 'BabyDogeCoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D23').Value = '0.0₆0712'
$ws.Range('E23').Value = '  +154.14%  '
$ws.Range('E24').Value = '  +9.50%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = "'65.96"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('D27').Value = "'648.80"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.39%  '
$ws.Range('D28').Value = "'1.18"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +17.52%  '
$ws.Range('D29').Value = "'8.49"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.09%  '
$ws.Range('D30').Value = '0.0₃0987'
$ws.Range('E30').Value = '  +5.46%  '
$ws.Range('D31').Value = '2.565.81'
$ws.Range('E31').Value = '  +2.07%  '
$ws.Range('D32').Value = "'1.46"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.49%  '
$ws.Range('D33').Value = "'8.20"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('E35').Value = '  +4.52%  '
$ws.Range('D36').Value = "'1.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.96%  '
$ws.Range('D37').Value = "'0.998"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = "'4.76"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.50%  '
$ws.Range('D39').Value = "'5.50"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.43%  '
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('D41').Value = "'153.20"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'18.74"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('E43').Value = '  +7.22%  '
$ws.Range('E44').Value = '  +4.67%  '
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = "'14.93"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +27.28%  '
$ws.Range('D48').Value = "'145.47"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').Value = "'3.63"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.29%  '
$ws.Range('D50').Value = "'20.63"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.43%  '
$ws.Range('D51').Value = "'0.605"
$ws.Range('D51').Style = 'Normal'
